$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 866.3333
$ws.Range("I18").Value = 866.3333
$ws.Range("K18").Value = 866.3333
$ws.Range("M18").Value = -582.3333

$ws.Range("H40").Value = 2603.5386
$ws.Range("J40").Value = 3471.2856
$ws.Range("L40").Value = 3471.2856
$ws.Range("N40").Value = -3821.2856

$ws.Range("H103").Value = 2166.5
$ws.Range("I103").Value = 2000
$ws.Range("K103").Value = 6000
$ws.Range("M103").Value = -5414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""

$ws.Range("H32").Value = 5587.089
$ws.Range("I32").Value = 5587.089
$ws.Range("K32").Value = 5587.089
$ws.Range("M32").Value = -5300.089

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = ""

$ws.Range("H61").Value = 2661.7
$ws.Range("I61").Value = 2661.7
$ws.Range("K61").Value = 2661.7
$ws.Range("M61").Value = -2449.7

$ws.Range("H63").Value = 4239.3477
$ws.Range("I63").Value = 2175.4
$ws.Range("K63").Value = 2175.4
$ws.Range("M63").Value = -1489.4

$ws.Range("H66").Value = 4239.3477
$ws.Range("I66").Value = 2175.4
$ws.Range("K66").Value = 10877
$ws.Range("M66").Value = -7445

$ws.Range("H74").Value = 3230.9443
$ws.Range("I74").Value = 3639.9285
$ws.Range("K74").Value = 3639.9285
$ws.Range("M74").Value = -2765.9285

$ws.Range("H77").Value = 3230.9443
$ws.Range("I77").Value = 3639.9285
$ws.Range("K77").Value = 18199.6425
$ws.Range("M77").Value = -13831.6425

$ws.Range("H102").Value = 1423.6666
$ws.Range("I102").Value = 1423.6666
$ws.Range("K102").Value = 1423.6666
$ws.Range("M102").Value = 198.3334

$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774

$ws.Range("H113").Value = 39500
$ws.Range("J113").Value = 39500
$ws.Range("L113").Value = 39500
$ws.Range("N113").Value = -48178

$ws.Range("H115").Value = 50684
$ws.Range("J115").Value = 50684
$ws.Range("L115").Value = 50684
$ws.Range("N115").Value = -53818

$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676

$ws.Range("H122").Value = 1604
$ws.Range("I122").Value = 1604
$ws.Range("K122").Value = 4812
$ws.Range("M122").Value = -2362

$ws.Range("H132").Value = 1636
$ws.Range("I132").Value = 1606.7778
$ws.Range("K132").Value = 4820.3334
$ws.Range("M132").Value = -2290.3334

$ws.Range("H136").Value = 2661.7
$ws.Range("I136").Value = 2661.7
$ws.Range("K136").Value = 7985.099999999999
$ws.Range("M136").Value = -5435.099999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""

$ws.Range("H36").Value = 1417.5
$ws.Range("I36").Value = 1223.3334
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 1223.3334
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = -689.3334
$ws.Range("N36").Value = -3068

$ws.Range("H64").Value = 2417.7144
$ws.Range("J64").Value = 2373.5386
$ws.Range("L64").Value = 2373.5386
$ws.Range("N64").Value = -2823.5386

$ws.Range("H67").Value = 2417.7144
$ws.Range("J67").Value = 2373.5386
$ws.Range("L67").Value = 2373.5386
$ws.Range("N67").Value = -3933.5386

$ws.Range("H103").Value = 17799.2
$ws.Range("J103").Value = 17799.2
$ws.Range("L103").Value = 17799.2
$ws.Range("N103").Value = -20143.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2273.5454
$ws.Range("I31").Value = 2332.6667
$ws.Range("J31").Value = 2251.375
$ws.Range("K31").Value = 2332.6667
$ws.Range("L31").Value = 2251.375
$ws.Range("M31").Value = -2037.6667
$ws.Range("N31").Value = -2841.375

$ws.Range("H34").Value = 2273.5454
$ws.Range("I34").Value = 2332.6667
$ws.Range("J34").Value = 2251.375
$ws.Range("K34").Value = 2332.6667
$ws.Range("L34").Value = 2251.375
$ws.Range("M34").Value = -2130.6667
$ws.Range("N34").Value = -2655.375

$ws.Range("H58").Value = 2500
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = ""

$ws.Range("H134").Value = 2872.0667
$ws.Range("I134").Value = 2898.5
$ws.Range("J134").Value = 2766.3333
$ws.Range("K134").Value = 8695.5
$ws.Range("L134").Value = 8298.999899999999
$ws.Range("M134").Value = -6160.5
$ws.Range("N134").Value = -13368.9999

$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 161.7
$ws.Range("J12").Value = 80.5
$ws.Range("L12").Value = 241.5
$ws.Range("N12").Value = -587.5

$ws.Range("H33").Value = 149.25
$ws.Range("J33").Value = 148
$ws.Range("L33").Value = 888
$ws.Range("N33").Value = -1454

$ws.Range("H55").Value = 2045.8
$ws.Range("J55").Value = 2739.5
$ws.Range("L55").Value = 8218.5
$ws.Range("N55").Value = -8572.5

$ws.Range("H133").Value = 10748.75
$ws.Range("I133").Value = 10748.75
$ws.Range("K133").Value = 32246.25
$ws.Range("M133").Value = -27186.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 32500
$ws.Range("J121").Value = 32500
$ws.Range("L121").Value = 32500
$ws.Range("N121").Value = -35994

$ws.Range("H122").Value = 2970.25
$ws.Range("I122").Value = 1612.2858
$ws.Range("J122").Value = 4871.4
$ws.Range("K122").Value = 4836.857400000001
$ws.Range("L122").Value = 14614.2
$ws.Range("M122").Value = -2386.857400000001
$ws.Range("N122").Value = -19514.2

$ws.Range("H123").Value = 55000
$ws.Range("J123").Value = 55000
$ws.Range("L123").Value = 55000
$ws.Range("N123").Value = -59900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2885.2856
$ws.Range("I22").Value = 3862.25
$ws.Range("K22").Value = 3862.25
$ws.Range("M22").Value = -3567.25

$ws.Range("H27").Value = 2885.2856
$ws.Range("I27").Value = 3862.25
$ws.Range("K27").Value = 3862.25
$ws.Range("M27").Value = -3755.25

$ws.Range("H46").Value = 3480.0908
$ws.Range("I46").Value = 2823.5
$ws.Range("K46").Value = 2823.5
$ws.Range("M46").Value = -2635.5

$ws.Range("H55").Value = 889.9231
$ws.Range("J55").Value = 1675
$ws.Range("L55").Value = 1675
$ws.Range("N55").Value = -2021

$ws.Range("H82").Value = 5042.7144
$ws.Range("I82").Value = 2659.8
$ws.Range("K82").Value = 2659.8
$ws.Range("M82").Value = -2298.8

$ws.Range("H85").Value = 5042.7144
$ws.Range("I85").Value = 2659.8
$ws.Range("K85").Value = 2659.8
$ws.Range("M85").Value = -1411.8

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""

$ws.Range("H122").Value = 3959.2
$ws.Range("I122").Value = 3959.2
$ws.Range("K122").Value = 11877.6
$ws.Range("M122").Value = -9427.599999999999

$ws.Range("H124").Value = 71996.60000000001
$ws.Range("J124").Value = 71996.60000000001
$ws.Range("L124").Value = 71996.60000000001
$ws.Range("N124").Value = -81816.60000000001

$ws.Range("H134").Value = 68499
$ws.Range("J134").Value = 68499
$ws.Range("L134").Value = 68499
$ws.Range("N134").Value = -78639

$ws.Range("H139").Value = 99998.836
$ws.Range("J139").Value = 99998.836
$ws.Range("L139").Value = 99998.836
$ws.Range("N139").Value = -110278.836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 250000740
$ws.Range("I100").Value = 250000740
$ws.Range("K100").Value = 500001480
$ws.Range("M100").Value = -500000939

$ws.Range("H122").Value = 4119.476
$ws.Range("I122").Value = 3970.3125
$ws.Range("J122").Value = 4596.8
$ws.Range("K122").Value = 11910.9375
$ws.Range("L122").Value = 13790.4
$ws.Range("M122").Value = -9460.9375
$ws.Range("N122").Value = -18690.4

$ws.Range("H132").Value = 3252
$ws.Range("I132").Value = 3314.6667
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 9944.000100000001
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -7414.000100000001
$ws.Range("N132").Value = -12560
